# Commit: "Add files via upload"
#
# The workbook's row 18 ("Age (in years)" section header) and row 28
# ("Education" section header) are reworded in all three languages so that
# each column reads as a "By ..." / "По ..." qualifier instead of a bare
# noun, e.g.
#   ky: "Жаш курагы (жылдарда)"  -> "Жаш курагы боюнча (жылдарда)"
#   ru: "Возраст (в годах)"      -> "По возрасту (в годах)"
#   en: "Age (in years) "        -> "By age (in years) "
#   ky: "Билими"                 -> "Билими боюнча"
#   ru: "Образование"            -> "По образованию"
#   en: "Education"              -> "By education"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18 - "age (in years)" section header, one cell per language
# (English updated first, then Russian, then Kyrgyz, to match the order
# the new labels were appended to the shared-string table.)
$ws.Cells.Item(18, 3).Value = "By age (in years) "

# Row 28 - "education" section header, English column
$ws.Cells.Item(28, 3).Value = "By education"

# Row 18 - Russian column
$ws.Cells.Item(18, 2).Value = "По возрасту (в годах)"

# Row 18 - Kyrgyz column
$ws.Cells.Item(18, 1).Value = "Жаш курагы боюнча (жылдарда)"

# Row 28 - Kyrgyz column
$ws.Cells.Item(28, 1).Value = "Билими боюнча"

# Row 28 - Russian column
$ws.Cells.Item(28, 2).Value = "По образованию"
